$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.712.69'
$ws.Range('E2').Value = '  -0.39%  '

$ws.Range('D3').Value = '1.633.13'
$ws.Range('E3').Value = '  -0.41%  '

$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').Value = '215.28'
$ws.Range('E5').Value = '  -0.22%  '

$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  -0.95%  '

$ws.Range('E7').Value = '  +0.19%  '

$ws.Range('E8').Value = '  -0.31%  '

$ws.Range('D9').Value = '0.0635'
$ws.Range('E9').Value = '  -1.33%  '

$ws.Range('D10').Value = '19.53'
$ws.Range('E10').Value = '  -4.52%  '

$ws.Range('E11').Value = '  +0.88%  '

$ws.Range('D12').Value = '4.23'
$ws.Range('E12').Value = '  -0.84%  '

$ws.Range('D13').Value = '1.858.74'
$ws.Range('E13').Value = '  -0.42%  '

$ws.Range('D14').Value = '1.631.98'
$ws.Range('E14').Value = '  -0.50%  '

$ws.Range('D15').Value = '0.554'
$ws.Range('E15').Value = '  -1.34%  '

$ws.Range('D16').Value = '0.0₃0765'
$ws.Range('E16').Value = '  -0.10%  '

$ws.Range('D17').Value = '62.66'
$ws.Range('E17').Value = '  -1.17%  '

$ws.Range('D18').Value = '25.749.15'
$ws.Range('E18').Value = '  -0.33%  '

$ws.Range('E19').Value = '  +0.14%  '

$ws.Range('D20').Value = '4.43'
$ws.Range('E20').Value = '  +0.89%  '

$ws.Range('D21').Value = '193.00'
$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').Value = '9.93'
$ws.Range('E22').Value = '  -0.09%  '

$ws.Range('D23').Value = '6.26'
$ws.Range('E23').Value = '  +1.87%  '

$ws.Range('E24').Value = '  +0.21%  '

$ws.Range('D25').Value = '1.84'
$ws.Range('E25').Value = '  +2.20%  '

$ws.Range('D26').Value = '140.20'
$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -1.85%  '

$ws.Range('D28').Value = '6.86'
$ws.Range('E28').Value = '  +0.45%  '

$ws.Range('D29').Value = '15.47'
$ws.Range('E29').Value = '  -0.93%  '

$ws.Range('E30').Value = '  -0.27%  '

$ws.Range('D31').Value = '0.0492'
$ws.Range('E31').Value = '  -1.05%  '

$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  +1.27%  '

$ws.Range('D33').Value = '3.24'
$ws.Range('E33').Value = '  -0.11%  '

$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +0.92%  '

$ws.Range('E35').Value = '  +0.32%  '

$ws.Range('E36').Value = '  -0.73%  '

$ws.Range('D37').Value = '0.548'
$ws.Range('E37').Value = '  -1.69%  '

$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.117.89'
$ws.Range('E38').Value = '  -1.34%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.52'
$ws.Range('E39').Value = '  -1.99%  '

$ws.Range('E40').Value = '  -0.87%  '

$ws.Range('E41').Value = '  +0.68%  '

$ws.Range('D42').Value = '5.55'
$ws.Range('E42').Value = '  +1.12%  '

$ws.Range('D43').Value = '99.64'
$ws.Range('E43').Value = '  +0.60%  '

$ws.Range('E44').Value = '  -0.15%  '

$ws.Range('D45').Value = '1.767.95'
$ws.Range('E45').Value = '  -0.47%  '

$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  +0.14%  '

$ws.Range('D47').Value = '55.06'
$ws.Range('E47').Value = '  -1.07%  '

$ws.Range('E48').Value = '  -2.24%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0501'
$ws.Range('E49').Value = '  -0.32%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.59'
$ws.Range('E50').Value = '  -2.61%  '

$ws.Range('E51').Value = '  +2.71%  '
